# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handback packages have been generated: it marks the
# Overview "Status" as handed back, records the handback datetimes, and
# fills in the "Latest Target File" / "Latest Handback File" (and, for
# de-de, "Latest Handback DateTime") columns on the language sheets, adding
# hyperlinks on the newly-populated "Latest Target File" cells.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns (zh-cn / de-de) for both rows ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn / de-de sheets: Status column (C) for both rows ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Widen the Overview Status columns (E, F) now that the text is longer
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# --- zh-cn sheet ---
# Row 2 (5b935408...md)
$zhcn.Range("I2").Value = "5b935408-90fc-4b55-a235-4bf052352988.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb780f092c522aa03698e494d6dfcb8754686810/e2e/5b935408-90fc-4b55-a235-4bf052352988.md", $null, $null, "5b935408-90fc-4b55-a235-4bf052352988.md")
$zhcn.Range("I2").Style = "Hyperlink"
$zhcn.Range("J2").Value = "5b935408-90fc-4b55-a235-4bf052352988.dd784210a8f4fe6da0a30a1a6105e46179c5fdee.zh-cn.xlf"

# Row 3 (fdbbaee0...md)
$zhcn.Range("I3").Value = "fdbbaee0-fc80-46b6-902e-aaf59d871475.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb780f092c522aa03698e494d6dfcb8754686810/e2e/fdbbaee0-fc80-46b6-902e-aaf59d871475.md", $null, $null, "fdbbaee0-fc80-46b6-902e-aaf59d871475.md")
$zhcn.Range("I3").Style = "Hyperlink"
$zhcn.Range("J3").Value = "fdbbaee0-fc80-46b6-902e-aaf59d871475.3cfbbee4c6899c7fb20f443789fa37db36059e6f.zh-cn.xlf"

# Handback datetime for zh-cn now in sync (was the "never handed back" placeholder)
$zhcn.Range("K2").Value = "2016-09-05 11:44:36"
$zhcn.Range("K3").Value = "2016-09-05 11:44:36"

# Widen the columns that now hold the longer handback filenames/status text
$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(9).ColumnWidth = 39.16666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.16666666666667

# --- de-de sheet ---
# Row 2 (5b935408...md)
$dede.Range("I2").Value = "5b935408-90fc-4b55-a235-4bf052352988.md"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb780f092c522aa03698e494d6dfcb8754686810/e2e/5b935408-90fc-4b55-a235-4bf052352988.md", $null, $null, "5b935408-90fc-4b55-a235-4bf052352988.md")
$dede.Range("I2").Style = "Hyperlink"
$dede.Range("J2").Value = "5b935408-90fc-4b55-a235-4bf052352988.dd784210a8f4fe6da0a30a1a6105e46179c5fdee.de-de.xlf"
$dede.Range("K2").Value = "2016-09-05 11:44:57"

# Row 3 (fdbbaee0...md)
$dede.Range("I3").Value = "fdbbaee0-fc80-46b6-902e-aaf59d871475.md"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb780f092c522aa03698e494d6dfcb8754686810/e2e/fdbbaee0-fc80-46b6-902e-aaf59d871475.md", $null, $null, "fdbbaee0-fc80-46b6-902e-aaf59d871475.md")
$dede.Range("I3").Style = "Hyperlink"
$dede.Range("J3").Value = "fdbbaee0-fc80-46b6-902e-aaf59d871475.3cfbbee4c6899c7fb20f443789fa37db36059e6f.de-de.xlf"
$dede.Range("K3").Value = "2016-09-05 11:44:57"

# Widen the columns that now hold the longer handback filenames/status text
$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(9).ColumnWidth = 39.16666666666667
$dede.Columns.Item(10).ColumnWidth = 39.16666666666667
